# Offer.xlsx update: merge the "System Architect" and "Software Developer"
# roles into a single "System Architect & Software Developer" row, drop the
# now-redundant standalone Software Developer row, and let every dependent
# formula / total follow the new numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Variable-cost table (A4:C11 -> A4:C10) -------------------------------

# Rename role in row 8 ("System Architect" -> merged role) and bump its
# hourly wage 90 -> 100; the existing formula in C8 (=B8*150) recalculates
# automatically.
$ws.Range("A8").Value = "System Architect & Software Developer"
$ws.Range("B8").Value = 100

# Remove the old standalone "Software Developer" row (old row 10); every row
# below it shifts up, the totals row lands on row 10, and table1's range
# shrinks from A4:C11 to A4:C10 automatically.
$ws.Range("A10").EntireRow.Delete()

# Re-apply the "Costs of 150 Hours" formula across the shrunk range so it is
# one shared formula group again (C5:C9).
$ws.Range("C5:C9").Formula = "=B5*150"

# The total-variable-costs formula referenced the now-removed row; point it
# at the remaining rows only.
$ws.Range("C10").Formula = "=C5+C6+C7+C8+C9"

# --- Fix-costs table: headcount per line item drops from 6 to 5 ----------

$ws.Range("G5").Value = 5
$ws.Range("G6").Value = 5

# --- Totals tables: swap hard-coded numbers for live formulas -------------

$ws.Range("K5").Formula = "=C10"
$ws.Range("K6").Formula = "=H8"
$ws.Range("N5").Formula = "=K7"

# --- Cosmetic / view updates (best effort) ---------------------------------

$ws.Columns("A:A").ColumnWidth = 32.17

$ws.Range("K7").Select()

$win = $excel.ActiveWindow
$win.Width = 23040
$win.Height = 9192
